$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.240179
$ws.Range("H2").Value = 3.720537
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.5292956666666667
$ws.Range("N2").Value = 1.587887
$ws.Range("O2").Value = 0.135651968140022
$ws.Range("P2").Value = 0.1356519681400219
$ws.Range("Q2").Value = 0.656421370591
$ws.Range("R2").Value = 5.907792335319
$ws.Range("S2").Value = 0.135651968140022
$ws.Range("T2").Value = 0.1356519681400219

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.240179
$ws.Range("H3").Value = 3.720537
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.362890666666667
$ws.Range("N3").Value = 4.088672000000001
$ws.Range("O3").Value = 0.3492921120199358
$ws.Range("P3").Value = 0.3492921120199358
$ws.Range("Q3").Value = 1.690228384096
$ws.Range("R3").Value = 15.212055456864
$ws.Range("S3").Value = 0.3492921120199358
$ws.Range("T3").Value = 0.3492921120199358

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.240179
$ws.Range("H4").Value = 3.720537
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.009678666666666
$ws.Range("N4").Value = 6.029036
$ws.Range("O4").Value = 0.5150559198400423
$ws.Range("P4").Value = 0.5150559198400423
$ws.Range("Q4").Value = 2.492361279148
$ws.Range("R4").Value = 22.431251512332
$ws.Range("S4").Value = 0.5150559198400423
$ws.Range("T4").Value = 0.5150559198400423
